# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures to the Famfrit_Profits workbook
# (chore: update Sheets via scheduled runner)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H11").Value = 136.9
$ws.Range("I11").Value = 136.9
$ws.Range("K11").Value = 136.9
$ws.Range("M11").Value = 3.099999999999994
$ws.Range("H15").Value = 582.8043
$ws.Range("I15").Value = 582.8043
$ws.Range("K15").Value = 1748.4129
$ws.Range("M15").Value = -1579.4129
$ws.Range("H28").Value = 7266.077
$ws.Range("I28").Value = 4606.5557
$ws.Range("J28").Value = 13250
$ws.Range("K28").Value = 4606.5557
$ws.Range("L28").Value = 13250
$ws.Range("M28").Value = -4121.5557
$ws.Range("N28").Value = -14220
$ws.Range("H32").Value = 3641.125
$ws.Range("I32").Value = 2556.8572
$ws.Range("J32").Value = 4484.4443
$ws.Range("K32").Value = 2556.8572
$ws.Range("L32").Value = 4484.4443
$ws.Range("M32").Value = -2230.8572
$ws.Range("N32").Value = -5136.4443
$ws.Range("H34").Value = 3866.6667
$ws.Range("I34").Value = 3866.6667
$ws.Range("K34").Value = 3866.6667
$ws.Range("M34").Value = -3663.6667
$ws.Range("H36").Value = 3866.6667
$ws.Range("I36").Value = 3866.6667
$ws.Range("K36").Value = 3866.6667
$ws.Range("M36").Value = -3151.6667
$ws.Range("H54").Value = 7000
$ws.Range("I54").Value = 7000
$ws.Range("K54").Value = 7000
$ws.Range("M54").Value = -6514
$ws.Range("H80").Value = 1993.3793
$ws.Range("I80").Value = 604.2
$ws.Range("K80").Value = 1812.6
$ws.Range("M80").Value = -814.6000000000001
$ws.Range("H83").Value = 1993.3793
$ws.Range("I83").Value = 604.2
$ws.Range("K83").Value = 5437.8
$ws.Range("M83").Value = -445.8000000000002
$ws.Range("H92").Value = 775.63635
$ws.Range("I92").Value = 673.2
$ws.Range("K92").Value = 673.2
$ws.Range("M92").Value = 574.8
$ws.Range("H96").Value = 19724.455
$ws.Range("I96").Value = 23329.945
$ws.Range("J96").Value = 3499.75
$ws.Range("K96").Value = 69989.83499999999
$ws.Range("L96").Value = 10499.25
$ws.Range("M96").Value = -68616.83499999999
$ws.Range("N96").Value = -13245.25
$ws.Range("H97").Value = 890
$ws.Range("J97").Value = 987.5
$ws.Range("L97").Value = 2962.5
$ws.Range("N97").Value = -3954.5
$ws.Range("H98").Value = 2866.4773
$ws.Range("I98").Value = 2681.2974
$ws.Range("J98").Value = 3845.2856
$ws.Range("K98").Value = 2681.2974
$ws.Range("L98").Value = 3845.2856
$ws.Range("M98").Value = -1183.2974
$ws.Range("N98").Value = -6841.2856
$ws.Range("H118").Value = 273.14285
$ws.Range("I118").Value = 273.14285
$ws.Range("K118").Value = 819.4285500000001
$ws.Range("M118").Value = 837.5714499999999
$ws.Range("H122").Value = 2866.4773
$ws.Range("I122").Value = 2681.2974
$ws.Range("J122").Value = 3845.2856
$ws.Range("K122").Value = 8043.8922
$ws.Range("L122").Value = 11535.8568
$ws.Range("M122").Value = -5593.8922
$ws.Range("N122").Value = -16435.8568
$ws.Range("H132").Value = 4046.65
$ws.Range("I132").Value = 3939.0571
$ws.Range("K132").Value = 11817.1713
$ws.Range("M132").Value = -9287.1713
$ws.Range("H138").Value = 4034.843
$ws.Range("J138").Value = 4324.049
$ws.Range("L138").Value = 12972.147
$ws.Range("N138").Value = -23252.147

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H74").Value = 32295632
$ws.Range("I74").Value = 38505852
$ws.Range("J74").Value = 2498.2
$ws.Range("K74").Value = 38505852
$ws.Range("L74").Value = 2498.2
$ws.Range("M74").Value = -38504978
$ws.Range("N74").Value = -4246.2
$ws.Range("H77").Value = 32295632
$ws.Range("I77").Value = 38505852
$ws.Range("J77").Value = 2498.2
$ws.Range("K77").Value = 192529260
$ws.Range("L77").Value = 12491
$ws.Range("M77").Value = -192524892
$ws.Range("N77").Value = -21227
$ws.Range("H102").Value = 252742
$ws.Range("I102").Value = 335331.16
$ws.Range("J102").Value = 4974.5
$ws.Range("K102").Value = 335331.16
$ws.Range("L102").Value = 4974.5
$ws.Range("M102").Value = -333709.16
$ws.Range("N102").Value = -8218.5
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = $null
$ws.Range("H122").Value = 10756907
$ws.Range("I122").Value = 2661.182
$ws.Range("K122").Value = 7983.545999999999
$ws.Range("M122").Value = -5533.545999999999

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H7").Value = 166668110
$ws.Range("I7").Value = 166668110
$ws.Range("K7").Value = 166668110
$ws.Range("M7").Value = -166667997
$ws.Range("H54").Value = 10343.375
$ws.Range("I54").Value = 1820.8572
$ws.Range("K54").Value = 1820.8572
$ws.Range("M54").Value = -1336.8572
$ws.Range("H94").Value = 544.59375
$ws.Range("I94").Value = 576.8570999999999
$ws.Range("J94").Value = 318.75
$ws.Range("K94").Value = 576.8570999999999
$ws.Range("L94").Value = 318.75
$ws.Range("M94").Value = -125.8570999999999
$ws.Range("N94").Value = -1220.75
$ws.Range("H107").Value = 2386.4
$ws.Range("J107").Value = 2775.7778
$ws.Range("L107").Value = 2775.7778
$ws.Range("N107").Value = -6615.7778
$ws.Range("H134").Value = 2886.037
$ws.Range("I134").Value = 2179.2
$ws.Range("K134").Value = 6537.599999999999
$ws.Range("M134").Value = -4002.599999999999

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H86").Value = 7504.231
$ws.Range("I86").Value = 7082.778
$ws.Range("K86").Value = 7082.778
$ws.Range("M86").Value = -5959.778
$ws.Range("H89").Value = 7504.231
$ws.Range("I89").Value = 7082.778
$ws.Range("K89").Value = 35413.89
$ws.Range("M89").Value = -29797.89
$ws.Range("H99").Value = 8235.788
$ws.Range("I99").Value = 1765.25
$ws.Range("J99").Value = 9128.275
$ws.Range("K99").Value = 1765.25
$ws.Range("L99").Value = 9128.275
$ws.Range("M99").Value = -267.25
$ws.Range("N99").Value = -12124.275
$ws.Range("H107").Value = 523.3077
$ws.Range("I107").Value = 482.5
$ws.Range("K107").Value = 482.5
$ws.Range("M107").Value = 1437.5
$ws.Range("H126").Value = 8235.788
$ws.Range("I126").Value = 1765.25
$ws.Range("J126").Value = 9128.275
$ws.Range("K126").Value = 5295.75
$ws.Range("L126").Value = 27384.825
$ws.Range("M126").Value = -2825.75
$ws.Range("N126").Value = -32324.825
$ws.Range("H132").Value = 108041.52
$ws.Range("I132").Value = 127730.75
$ws.Range("K132").Value = 383192.25
$ws.Range("M132").Value = -380662.25

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H98").Value = 753
$ws.Range("J98").Value = 750
$ws.Range("L98").Value = 2250
$ws.Range("N98").Value = -5246

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H46").Value = 3500
$ws.Range("I46").Value = 3500
$ws.Range("K46").Value = 3500
$ws.Range("M46").Value = -3344
$ws.Range("H70").Value = 232868.11
$ws.Range("I70").Value = 343468.84
$ws.Range("K70").Value = 343468.84
$ws.Range("M70").Value = -343198.84
$ws.Range("H73").Value = 232868.11
$ws.Range("I73").Value = 343468.84
$ws.Range("K73").Value = 343468.84
$ws.Range("M73").Value = -342532.84
$ws.Range("H93").Value = 33163.332
$ws.Range("J93").Value = 32490
$ws.Range("L93").Value = 32490
$ws.Range("N93").Value = -36234
$ws.Range("H97").Value = 998.2857
$ws.Range("I97").Value = 698
$ws.Range("K97").Value = 698
$ws.Range("M97").Value = -202
$ws.Range("H107").Value = 806
$ws.Range("I107").Value = 871.1429000000001
$ws.Range("K107").Value = 871.1429000000001
$ws.Range("M107").Value = 1048.8571

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H53").Value = 15023
$ws.Range("I53").Value = 46
$ws.Range("J53").Value = 30000
$ws.Range("K53").Value = 46
$ws.Range("L53").Value = 30000
$ws.Range("M53").Value = 472
$ws.Range("N53").Value = -31036
$ws.Range("H68").Value = 4496.5
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").Value = $null
$ws.Range("H71").Value = 4496.5
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").Value = $null

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H37").Value = 21249.5
$ws.Range("J37").Value = 21249.5
$ws.Range("L37").Value = 21249.5
$ws.Range("N37").Value = -21655.5
$ws.Range("H107").Value = 1572.3214
$ws.Range("I107").Value = 1242.8667
$ws.Range("K107").Value = 3728.6001
$ws.Range("M107").Value = -1808.6001
$ws.Range("H113").Value = 1147.1666
$ws.Range("I113").Value = 1264.1428
$ws.Range("J113").Value = 983.4
$ws.Range("K113").Value = 3792.4284
$ws.Range("L113").Value = 2950.2
$ws.Range("M113").Value = -1622.4284
$ws.Range("N113").Value = -7290.2
$ws.Range("H122").Value = 5265578
$ws.Range("I122").Value = 1862.2646
$ws.Range("J122").Value = 50007164
$ws.Range("K122").Value = 5586.793799999999
$ws.Range("L122").Value = 150021492
$ws.Range("M122").Value = -3136.793799999999
$ws.Range("N122").Value = -150026392

Write-Host "Updated profit figures across ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR sheets."